# Append: 2025-10-26 06:29 JST
# Updates the "ランサーズ" (Lancers) listing sheet: refreshes the two newest
# rows with new scrape data and drops the now-stale rows 4-14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- Update row 2 --------------------------------------------------------
$ws.Range("A2").Value = "2025-10-26 06:29:44"
$ws.Range("B2").Value = "【機密性の高いノウハウを含む】サーバーレスAI分析システム構築(MVP開発と拡張性確保)"
$ws.Range("D2").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5420678"
$ws.Range("G2").Value = 383
$ws.Range("H2").Value = "🔥AI,Ai ◆開発"

# --- Update row 3 --------------------------------------------------------
$ws.Range("A3").Value = "2025-10-26 06:29:44"
$ws.Range("B3").Value = "【Power Automate for Desktop】販売管理システムへExcelから自動入力"
$ws.Range("D3").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5407216"
$ws.Range("G3").Value = 48
$ws.Range("H3").Value = "◇管理"

# --- Drop the old rows 4-14 (only the two refreshed listings remain) -----
$ws.Range("A4:H14").EntireRow.Delete()

# --- Rebuild hyperlinks for the surviving rows ----------------------------
# Row deletion leaves stale hyperlink entries pointing past the new used
# range, and the F2/F3 URLs themselves changed - clear everything and
# re-add the two links that should remain.
$ws.Cells.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5420678")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5407216")

# --- Column width tweaks ---------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 28
$ws.Columns.Item(8).ColumnWidth = 12
